$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update video filenames: downloaded as .avi instead of .mp4 (some renamed with resolution)
$ws.Range("E5").Value = "Inscapes_02.avi"
$ws.Range("E6").Value = "The_Present_720x480.avi"
$ws.Range("E7").Value = "The_Present_720x480.avi"
$ws.Range("E8").Value = "Despicable_Me_720x480_English.avi"
$ws.Range("E9").Value = "Despicable_Me_720x480_English.avi"
$ws.Range("E10").Value = "Despicable_Me_720x480_Hungarian.avi"
$ws.Range("E11").Value = "Despicable_Me_720x480_Hungarian.avi"

# Move active selection from F17 to F15
$ws.Range("F15").Select()
